# CotisationCnas.docx : import rappel fait from xlsx
#
# The "rappel" table refers to Avril instead of Mars, and the
# amounts (number of months-equivalent/day rate, assiette total,
# 5% contribution) are bumped to the new, larger figures. The
# "sum in words" line at the bottom is updated to match.
#
# Each Find/Execute below uses Replace:=2 (wdReplaceAll) so every
# matching occurrence in the document (the values appear twice in
# the table: once in the data row, once in the "Total" row) gets
# updated in a single call.

$d = $word.ActiveDocument

# 1) Month label: Mars -> Avril (1 occurrence)
$d.Content.Find.Execute("Mars", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Avril", 2)

# 2) Rate/day-count column: 3868 -> 4053 (2 occurrences: data row + total row)
$d.Content.Find.Execute("3868", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "4053", 2)

# 3) Assiette total column: 69 624 000,00 -> 72 954 000,00 (2 occurrences)
$d.Content.Find.Execute("69 624 000,00", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "72 954 000,00", 2)

# 4) 5% contribution amount column: 3 481 200,00 -> 3 647 700,00 (2 occurrences)
$d.Content.Find.Execute("3 481 200,00", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "3 647 700,00", 2)

# 5) Amount spelled out in French words at the bottom of the document
$d.Content.Find.Execute("TROIS MILLIONS QUATRE CENT QUATRE-VINGT-UN MILLE DEUX CENTS  ", `
                         $true, $false, $false, $false, $false, `
                         $true, 1, $false, "TROIS MILLIONS SIX CENT QUARANTE-SEPT MILLE SEPT CENTS  ", 2)
